$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report date label (stored as plain text, not a date value)
$ws.Range("B1").Value = "25/03/2023"

# Update activity totals and hour counts
$ws.Range("B2").Value = 1666
$ws.Range("C2").Value = 14

$ws.Range("B3").Value = 107
$ws.Range("C3").Value = 14

$ws.Range("B4").Value = 63
$ws.Range("C4").Value = 14

$ws.Range("B5").Value = 821
$ws.Range("C5").Value = 14

$ws.Range("B6").Value = 564
$ws.Range("C6").Value = 14

$ws.Range("B7").Value = 174
$ws.Range("C7").Value = 14

$ws.Range("B8").Value = 120
$ws.Range("C8").Value = 14

$ws.Range("B9").Value = 778
$ws.Range("C9").Value = 14

$ws.Range("B10").Value = 73
$ws.Range("C10").Value = 14

$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 14

$ws.Range("B12").Value = 52
$ws.Range("C12").Value = 14
